$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 12
$ws.Range("H12").Value = 183.75
$ws.Range("I12").Value = 169.2
$ws.Range("K12").Value = 169.2
$ws.Range("M12").Value = 0.8000000000000114
# Row 19
$ws.Range("H19").Value = 13822.8125
$ws.Range("I19").Value = 2807
$ws.Range("J19").Value = 20432.3
$ws.Range("K19").Value = 2807
$ws.Range("L19").Value = 20432.3
$ws.Range("M19").Value = -2632
$ws.Range("N19").Value = -20782.3
# Row 46
$ws.Range("H46").Value = 11187.5
$ws.Range("I46").Value = 2375
$ws.Range("K46").Value = 7125
$ws.Range("M46").Value = -7006
# Row 60
$ws.Range("H60").Value = 11187.5
$ws.Range("I60").Value = 2375
$ws.Range("K60").Value = 7125
$ws.Range("M60").Value = -6641
# Row 100
$ws.Range("H100").Value = 7070.9165
$ws.Range("I100").Value = 8550.143
$ws.Range("K100").Value = 8550.143
$ws.Range("M100").Value = -8009.143
# Row 106
$ws.Range("H106").Value = 9233
$ws.Range("I106").Value = 2585.2856
$ws.Range("K106").Value = 2585.2856
$ws.Range("M106").Value = -1954.2856
# Row 107
$ws.Range("H107").Value = 462.5
$ws.Range("I107").Value = 450
$ws.Range("J107").Value = 500
$ws.Range("K107").Value = 450
$ws.Range("L107").Value = 500
$ws.Range("M107").Value = 1470
$ws.Range("N107").Value = -4340
# Row 137
$ws.Range("H137").Value = 2618
$ws.Range("I137").Value = 1552.75
$ws.Range("J137").Value = 3470.2
$ws.Range("K137").Value = 4658.25
$ws.Range("L137").Value = 10410.6
$ws.Range("M137").Value = -2108.25
$ws.Range("N137").Value = -15510.6
# Row 138
$ws.Range("H138").Value = 2003
$ws.Range("I138").Value = 1305.4584
$ws.Range("K138").Value = 3916.3752
$ws.Range("M138").Value = 1223.6248
# Row 141
$ws.Range("H141").Value = 3172.1738
$ws.Range("I141").Value = 1061.6316
$ws.Range("K141").Value = 3184.8948
$ws.Range("M141").Value = 1995.1052

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 2113.2222
$ws.Range("I45").Value = 2002.375
$ws.Range("K45").Value = 2002.375
$ws.Range("M45").Value = -1625.375
# Row 61
$ws.Range("H61").Value = 15876395
$ws.Range("I61").Value = 20836394
$ws.Range("K61").Value = 20836394
$ws.Range("M61").Value = -20836182
# Row 74
$ws.Range("H74").Value = 2861.25
$ws.Range("I74").Value = 2897.9092
$ws.Range("K74").Value = 2897.9092
$ws.Range("M74").Value = -2023.9092
# Row 77
$ws.Range("H77").Value = 2861.25
$ws.Range("I77").Value = 2897.9092
$ws.Range("K77").Value = 14489.546
$ws.Range("M77").Value = -10121.546
# Row 132
$ws.Range("H132").Value = 25001684
$ws.Range("I132").Value = 26317424
$ws.Range("K132").Value = 78952272
$ws.Range("M132").Value = -78949742
# Row 136
$ws.Range("H136").Value = 15876395
$ws.Range("I136").Value = 20836394
$ws.Range("K136").Value = 62509182
$ws.Range("M136").Value = -62506632

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 11414.967
$ws.Range("I20").Value = 14388.3
$ws.Range("J20").Value = 5468.3
$ws.Range("K20").Value = 14388.3
$ws.Range("L20").Value = 5468.3
$ws.Range("M20").Value = -14141.3
$ws.Range("N20").Value = -5962.3
# Row 80
$ws.Range("H80").Value = 31251098
$ws.Range("J80").Value = 41667590
$ws.Range("L80").Value = 41667590
$ws.Range("N80").Value = -41669586
# Row 83
$ws.Range("H83").Value = 31251098
$ws.Range("J83").Value = 41667590
$ws.Range("L83").Value = 208337950
$ws.Range("N83").Value = -208347934
# Row 86
$ws.Range("H86").Value = 2020.2759
$ws.Range("J86").Value = 2987.125
$ws.Range("L86").Value = 2987.125
$ws.Range("N86").Value = -5233.125
# Row 89
$ws.Range("H89").Value = 2020.2759
$ws.Range("J89").Value = 2987.125
$ws.Range("L89").Value = 14935.625
$ws.Range("N89").Value = -26167.625
# Row 99
$ws.Range("H99").Value = 1625.4667
$ws.Range("I99").Value = 1664.25
$ws.Range("K99").Value = 1664.25
$ws.Range("M99").Value = -166.25
# Row 134
$ws.Range("H134").Value = 1044.5
$ws.Range("I134").Value = 929.92
$ws.Range("K134").Value = 2789.76
$ws.Range("M134").Value = -254.7599999999998
# Row 138
$ws.Range("H138").Value = 60000
$ws.Range("J138").Value = 60000
$ws.Range("L138").Value = 60000
$ws.Range("N138").Value = -70280

$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 357.2
$ws.Range("I22").Value = 262
$ws.Range("J22").Value = 500
$ws.Range("K22").Value = 262
$ws.Range("L22").Value = 500
$ws.Range("M22").Value = 88
$ws.Range("N22").Value = -1200
# Row 31
$ws.Range("H31").Value = 2855.889
$ws.Range("I31").Value = 1767.5927
$ws.Range("J31").Value = 3508.8667
$ws.Range("K31").Value = 1767.5927
$ws.Range("L31").Value = 3508.8667
$ws.Range("M31").Value = -1472.5927
$ws.Range("N31").Value = -4098.8667
# Row 34
$ws.Range("H34").Value = 2855.889
$ws.Range("I34").Value = 1767.5927
$ws.Range("J34").Value = 3508.8667
$ws.Range("K34").Value = 1767.5927
$ws.Range("L34").Value = 3508.8667
$ws.Range("M34").Value = -1565.5927
$ws.Range("N34").Value = -3912.8667
# Row 58
$ws.Range("H58").Value = 2738.96
$ws.Range("J58").Value = 4303.9165
$ws.Range("L58").Value = 4303.9165
$ws.Range("N58").Value = -4709.9165
# Row 62
$ws.Range("H62").Value = 71435280
$ws.Range("I62").Value = 6789.3
$ws.Range("K62").Value = 6789.3
$ws.Range("M62").Value = -6165.3
# Row 65
$ws.Range("H65").Value = 71435280
$ws.Range("I65").Value = 6789.3
$ws.Range("K65").Value = 33946.5
$ws.Range("M65").Value = -30826.5
# Row 122
$ws.Range("H122").Value = 1991.762
$ws.Range("I122").Value = 2165.4707
$ws.Range("K122").Value = 6496.4121
$ws.Range("M122").Value = -4046.4121
# Row 132
$ws.Range("H132").Value = 1786.0646
$ws.Range("I132").Value = 1816.4333
$ws.Range("J132").Value = 875
$ws.Range("K132").Value = 5449.2999
$ws.Range("L132").Value = 2625
$ws.Range("M132").Value = -2919.2999
$ws.Range("N132").Value = -7685
# Row 134
$ws.Range("H134").Value = 2162.32
$ws.Range("I134").Value = 1239.0588
$ws.Range("J134").Value = 4124.25
$ws.Range("K134").Value = 3717.1764
$ws.Range("L134").Value = 12372.75
$ws.Range("M134").Value = -1182.1764
$ws.Range("N134").Value = -17442.75
# Row 136
$ws.Range("H136").Value = 2738.96
$ws.Range("J136").Value = 4303.9165
$ws.Range("L136").Value = 12911.7495
$ws.Range("N136").Value = -18011.7495

$ws = $wb.Worksheets.Item("CUL")
# Row 8
$ws.Range("H8").Value = 37710
$ws.Range("I8").Value = 37710
$ws.Range("K8").Value = 113130
$ws.Range("M8").Value = -112991
# Row 12
$ws.Range("H12").Value = 138.94118
$ws.Range("I12").Value = 188
$ws.Range("K12").Value = 564
$ws.Range("M12").Value = -391
# Row 132
$ws.Range("H132").Value = 1797.2
$ws.Range("I132").Value = 1300
$ws.Range("J132").Value = 1921.5
$ws.Range("K132").Value = 11700
$ws.Range("L132").Value = 17293.5
$ws.Range("M132").Value = -9170
$ws.Range("N132").Value = -22353.5

$ws = $wb.Worksheets.Item("GSM")
# Row 46
$ws.Range("H46").Value = 41555.445
$ws.Range("J46").Value = 41555.445
$ws.Range("L46").Value = 41555.445
$ws.Range("N46").Value = -41867.445
# Row 132
$ws.Range("H132").Value = 4738
$ws.Range("I132").Value = 4212.615
$ws.Range("J132").Value = 5713.7144
$ws.Range("K132").Value = 12637.845
$ws.Range("L132").Value = 17141.1432
$ws.Range("M132").Value = -10107.845
$ws.Range("N132").Value = -22201.1432
# Row 140
$ws.Range("H140").Value = 0
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 0
$ws.Range("L140").ClearContents()
$ws.Range("M140").ClearContents()
$ws.Range("N140").Value = 0

$ws = $wb.Worksheets.Item("LTW")
# Row 68
$ws.Range("H68").Value = 7425.25
$ws.Range("J68").Value = 15534.833
$ws.Range("L68").Value = 15534.833
$ws.Range("N68").Value = -17032.833
# Row 71
$ws.Range("H71").Value = 7425.25
$ws.Range("J71").Value = 15534.833
$ws.Range("L71").Value = 77674.16500000001
$ws.Range("N71").Value = -85162.16500000001
# Row 99
$ws.Range("H99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("L99").ClearContents()
$ws.Range("N99").Value = 0
# Row 132
$ws.Range("H132").Value = 6740.4
$ws.Range("I132").Value = 3536.3635
$ws.Range("J132").Value = 10656.444
$ws.Range("K132").Value = 10609.0905
$ws.Range("L132").Value = 31969.332
$ws.Range("M132").Value = -8079.0905
$ws.Range("N132").Value = -37029.33199999999
# Row 134
$ws.Range("H134").Value = 41011.8
$ws.Range("J134").Value = 41011.8
$ws.Range("L134").Value = 41011.8
$ws.Range("N134").Value = -51151.8
# Row 136
$ws.Range("H136").Value = 2287.3125
$ws.Range("I136").Value = 2153.7693
$ws.Range("J136").Value = 2866
$ws.Range("K136").Value = 6461.3079
$ws.Range("L136").Value = 8598
$ws.Range("M136").Value = -3911.3079
$ws.Range("N136").Value = -13698

$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 707.5
$ws.Range("I107").Value = 717.0769
$ws.Range("K107").Value = 2151.2307
$ws.Range("M107").Value = -231.2307000000001
# Row 136
$ws.Range("H136").Value = 3160.7083
$ws.Range("I136").Value = 1214.7222
$ws.Range("J136").Value = 8998.666999999999
$ws.Range("K136").Value = 3644.1666
$ws.Range("L136").Value = 26996.001
$ws.Range("M136").Value = -1094.1666
